# Automatic update of files.
# Re-sync species-observation rows 5-17 on the active sheet:
#  - rows are re-ordered (each row's data now comes from a different
#    source row after the refresh)
#  - coordinate columns Q (Ost) / R (Nord) are rounded to whole metres
#  - the per-row time-of-day columns Z (Starttid) / AB (Sluttid) are
#    dropped entirely
#  - the "Ålder-Stadium" column K keeps only a single value (row 15)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetRows = @(
    [PSCustomObject]@{ Row=5; A=111943981; B=96253; D="LC"; E=504; F="Guckusko"; G="Cypripedium calceolus"; H="L."; K=""; Q=682877; R=6694410 },
    [PSCustomObject]@{ Row=6; A=111943998; B=98535; D="LC"; E=222498; F="Blåsippa"; G="Hepatica nobilis"; H="Schreb."; K=""; Q=682757; R=6694406 },
    [PSCustomObject]@{ Row=7; A=111943999; B=99413; D="LC"; E=221235; F="Vårärt"; G="Lathyrus vernus"; H="(L.) Bernh."; K=""; Q=682757; R=6694406 },
    [PSCustomObject]@{ Row=8; A=111943980; B=89183; D="LC"; E=3215; F="Rödgul trumpetsvamp"; G="Craterellus lutescens"; H="(Fr.) Fr."; K=""; Q=682877; R=6694410 },
    [PSCustomObject]@{ Row=9; A=111943984; B=99413; D="LC"; E=221235; F="Vårärt"; G="Lathyrus vernus"; H="(L.) Bernh."; K=""; Q=682929; R=6694685 },
    [PSCustomObject]@{ Row=10; A=111943995; B=88899; D="NT"; E=3286; F="Flattoppad klubbsvamp"; G="Clavariadelphus truncatus"; H="(Quél.) Donk"; K=""; Q=682779; R=6694551 },
    [PSCustomObject]@{ Row=11; A=111943992; B=89183; D="LC"; E=3215; F="Rödgul trumpetsvamp"; G="Craterellus lutescens"; H="(Fr.) Fr."; K=""; Q=682867; R=6694644 },
    [PSCustomObject]@{ Row=12; A=111943979; B=96253; D="LC"; E=504; F="Guckusko"; G="Cypripedium calceolus"; H="L."; K=""; Q=682879; R=6694407 },
    [PSCustomObject]@{ Row=13; A=111943988; B=107033; D="NT"; E=220320; F="Ängsskära"; G="Serratula tinctoria"; H="L."; K=""; Q=682930; R=6694720 },
    [PSCustomObject]@{ Row=14; A=111943996; B=90332; D="LC"; E=4769; F="Svavelriska"; G="Lactarius scrobiculatus"; H="(Scop.:Fr.) Fr."; K=""; Q=682785; R=6694547 },
    [PSCustomObject]@{ Row=15; A=111943997; B=96326; D="LC"; E=219798; F="Skogsknipprot"; G="Epipactis helleborine"; H="(L.) Crantz"; K="i frukt"; Q=682781; R=6694488 },
    [PSCustomObject]@{ Row=16; A=111943990; B=101703; D="LC"; E=222412; F="Tibast"; G="Daphne mezereum"; H="L."; K=""; Q=682930; R=6694720 },
    [PSCustomObject]@{ Row=17; A=111943983; B=90678; D="LC"; E=4366; F="Skarp dropptaggsvamp"; G="Hydnellum peckii"; H="Banker"; K=""; Q=682871; R=6694481 }
)

foreach ($r in $targetRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value  = $r.A   # A: Id
    $ws.Cells.Item($row, 2).Value  = $r.B   # B: Taxonsorteringsordning
    $ws.Cells.Item($row, 4).Value  = $r.D   # D: Rödlistade
    $ws.Cells.Item($row, 5).Value  = $r.E   # E: TaxonId
    $ws.Cells.Item($row, 6).Value  = $r.F   # F: Artnamn
    $ws.Cells.Item($row, 7).Value  = $r.G   # G: Vetenskapligt namn
    $ws.Cells.Item($row, 8).Value  = $r.H   # H: Auktor

    if ($r.K -ne "") {
        $ws.Cells.Item($row, 11).Value = $r.K   # K: Ålder-Stadium
    } else {
        $ws.Cells.Item($row, 11).Value = $null
    }

    $ws.Cells.Item($row, 17).Value = $r.Q   # Q: Ost (rounded)
    $ws.Cells.Item($row, 18).Value = $r.R   # R: Nord (rounded)

    $ws.Cells.Item($row, 26).Value = $null  # Z: Starttid - removed
    $ws.Cells.Item($row, 28).Value = $null  # AB: Sluttid - removed
}
